$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 293.72
$ws.Range("I15").Value = 293.72
$ws.Range("K15").Value = 881.1600000000001
$ws.Range("M15").Value = -712.1600000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 10817615
$ws.Range("I135").Value = 328.33334
$ws.Range("J135").Value = 33128268
$ws.Range("K135").Value = 2955.00006
$ws.Range("L135").Value = 298154412
$ws.Range("M135").Value = -420.0000600000003
$ws.Range("N135").Value = -298159482

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 22728586
$ws.Range("I137").Value = 1197.258
$ws.Range("K137").Value = 3591.774
$ws.Range("M137").Value = -1041.774

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2150.4062
$ws.Range("I138").Value = 1617.0244
$ws.Range("J138").Value = 3101.2173
$ws.Range("K138").Value = 4851.0732
$ws.Range("L138").Value = 9303.651899999999
$ws.Range("M138").Value = 288.9268000000002
$ws.Range("N138").Value = -19583.6519

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4390.5244
$ws.Range("I32").Value = 4480.887
$ws.Range("K32").Value = 4480.887
$ws.Range("M32").Value = -4193.887

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5849204.5
$ws.Range("I61").Value = 6945636.5
$ws.Range("J61").Value = 1566.6666
$ws.Range("K61").Value = 6945636.5
$ws.Range("L61").Value = 1566.6666
$ws.Range("M61").Value = -6945424.5
$ws.Range("N61").Value = -1990.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 23152762
$ws.Range("I74").Value = 32052102
$ws.Range("J74").Value = 14474.5
$ws.Range("K74").Value = 32052102
$ws.Range("L74").Value = 14474.5
$ws.Range("M74").Value = -32051228
$ws.Range("N74").Value = -16222.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 23152762
$ws.Range("I77").Value = 32052102
$ws.Range("J77").Value = 14474.5
$ws.Range("K77").Value = 160260510
$ws.Range("L77").Value = 72372.5
$ws.Range("M77").Value = -160256142
$ws.Range("N77").Value = -81108.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2350
$ws.Range("I102").Value = 2270
$ws.Range("J102").Value = 2590
$ws.Range("K102").Value = 2270
$ws.Range("L102").Value = 2590
$ws.Range("M102").Value = -648
$ws.Range("N102").Value = -5834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 911792.7
$ws.Range("I132").Value = 1381515.5
$ws.Range("J132").Value = 73001.92999999999
$ws.Range("K132").Value = 4144546.5
$ws.Range("L132").Value = 219005.79
$ws.Range("M132").Value = -4142016.5
$ws.Range("N132").Value = -224065.79

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5849204.5
$ws.Range("I136").Value = 6945636.5
$ws.Range("J136").Value = 1566.6666
$ws.Range("K136").Value = 20836909.5
$ws.Range("L136").Value = 4699.9998
$ws.Range("M136").Value = -20834359.5
$ws.Range("N136").Value = -9799.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4281575
$ws.Range("I134").Value = 5295374
$ws.Range("K134").Value = 15886122
$ws.Range("M134").Value = -15883587

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2275.3403
$ws.Range("I31").Value = 1101.2424
$ws.Range("J31").Value = 5042.857
$ws.Range("K31").Value = 1101.2424
$ws.Range("L31").Value = 5042.857
$ws.Range("M31").Value = -806.2424000000001
$ws.Range("N31").Value = -5632.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2275.3403
$ws.Range("I34").Value = 1101.2424
$ws.Range("J34").Value = 5042.857
$ws.Range("K34").Value = 1101.2424
$ws.Range("L34").Value = 5042.857
$ws.Range("M34").Value = -899.2424000000001
$ws.Range("N34").Value = -5446.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1854.2778
$ws.Range("I58").Value = 857.4681
$ws.Range("J58").Value = 3728.28
$ws.Range("K58").Value = 857.4681
$ws.Range("L58").Value = 3728.28
$ws.Range("M58").Value = -654.4681
$ws.Range("N58").Value = -4134.280000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1813.5834
$ws.Range("I132").Value = 1774.275
$ws.Range("J132").Value = 2010.125
$ws.Range("K132").Value = 5322.825000000001
$ws.Range("L132").Value = 6030.375
$ws.Range("M132").Value = -2792.825000000001
$ws.Range("N132").Value = -11090.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1243.8914
$ws.Range("I134").Value = 1374.2
$ws.Range("J134").Value = 829.2727
$ws.Range("K134").Value = 4122.6
$ws.Range("L134").Value = 2487.8181
$ws.Range("M134").Value = -1587.6
$ws.Range("N134").Value = -7557.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1854.2778
$ws.Range("I136").Value = 857.4681
$ws.Range("J136").Value = 3728.28
$ws.Range("K136").Value = 2572.4043
$ws.Range("L136").Value = 11184.84
$ws.Range("M136").Value = -22.40430000000015
$ws.Range("N136").Value = -16284.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9616373
$ws.Range("I5").Value = 1233.5
$ws.Range("J5").Value = 12500915
$ws.Range("K5").Value = 3700.5
$ws.Range("L5").Value = 37502745
$ws.Range("M5").Value = -3588.5
$ws.Range("N5").Value = -37502969

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5316.5415
$ws.Range("I131").Value = 5504.1665
$ws.Range("J131").Value = 5128.9165
$ws.Range("K131").Value = 16512.4995
$ws.Range("L131").Value = 15386.7495
$ws.Range("M131").Value = -11472.4995
$ws.Range("N131").Value = -25466.7495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 9616373
$ws.Range("I135").Value = 1233.5
$ws.Range("J135").Value = 12500915
$ws.Range("K135").Value = 11101.5
$ws.Range("L135").Value = 112508235
$ws.Range("M135").Value = -8566.5
$ws.Range("N135").Value = -112513305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4357.4546
$ws.Range("I70").Value = 4116.5
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 4116.5
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -3846.5
$ws.Range("N70").Value = -5540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4357.4546
$ws.Range("I73").Value = 4116.5
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 4116.5
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -3180.5
$ws.Range("N73").Value = -6872

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3501.9424
$ws.Range("I122").Value = 2691.973
$ws.Range("J122").Value = 5499.8667
$ws.Range("K122").Value = 8075.919
$ws.Range("L122").Value = 16499.6001
$ws.Range("M122").Value = -5625.919
$ws.Range("N122").Value = -21399.6001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1698.525
$ws.Range("I132").Value = 1736.8064
$ws.Range("J132").Value = 1566.6666
$ws.Range("K132").Value = 5210.4192
$ws.Range("L132").Value = 4699.9998
$ws.Range("M132").Value = -2680.4192
$ws.Range("N132").Value = -9759.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2133.4856
$ws.Range("I16").Value = 1706.8387
$ws.Range("J16").Value = 5440
$ws.Range("K16").Value = 1706.8387
$ws.Range("L16").Value = 5440
$ws.Range("M16").Value = -1536.8387
$ws.Range("N16").Value = -5780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 44156.54
$ws.Range("I132").Value = 44156.54
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 132469.62
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -129939.62
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 29700
$ws.Range("J134").Value = 29700
$ws.Range("L134").Value = 29700
$ws.Range("N134").Value = -39840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 36514.266
$ws.Range("J123").Value = 36514.266
$ws.Range("L123").Value = 36514.266
$ws.Range("N123").Value = -46314.266

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3929.261
$ws.Range("I132").Value = 4593.6313
$ws.Range("J132").Value = 773.5
$ws.Range("K132").Value = 13780.8939
$ws.Range("L132").Value = 2320.5
$ws.Range("M132").Value = -11250.8939
$ws.Range("N132").Value = -7380.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6522.4116
$ws.Range("I136").Value = 7259.433
$ws.Range("K136").Value = 21778.299
$ws.Range("M136").Value = -19228.299
